# Double up the apostrophes in the "name" column (column B) for the four
# episode names that contain an apostrophe, matching the author's commit
# (looks like an escaping pass, e.g. for SQL insertion, was applied to a
# handful of rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = "Ransei''s Legendary Beauty"
$ws.Range("B28").Value = "The Visionary''s Quest"
$ws.Range("B31").Value = "The Dragon''s Dream"
$ws.Range("B38").Value = "The Free Spirit''s Path"
